$p = $ppt.ActivePresentation

# --- Update cached date fields (handout master + notes master) ---
# Handout master footer date: 25.06.2022 -> 02.07.2022
$hm = $p.HandoutMaster
$hmDate = $hm.HeadersFooters.DateAndTime
$hmDate.Text = "02.07.2022"

# Notes master footer date: 25/06/2022 -> 02/07/2022
$nm = $p.NotesMaster
$nmDate = $nm.HeadersFooters.DateAndTime
$nmDate.Text = "02/07/2022"

# --- Slide 3: "Group work (7 groups ..." -> "Group work (8 groups ..." ---
$s3 = $p.Slides.Item(3)
$s3Shape = $s3.Shapes.Item(5)
$s3OrigHeight = $s3Shape.Height
$s3OrigWidth = $s3Shape.Width
$s3Text = $s3Shape.TextFrame.TextRange
$idx = $s3Text.Text.IndexOf(" (7 ")
$s3Sub = $s3Text.Characters($idx + 1, 4)
$s3Sub.Text = " (8 "
# editing the run can nudge the auto-fit shape size; put it back exactly
$s3Shape.Height = $s3OrigHeight
$s3Shape.Width = $s3OrigWidth

# --- Slide 15: " Other?" -> " (Topic 8) Scheduling the Trading Bot" ---
$s15 = $p.Slides.Item(15)
$s15Shape = $s15.Shapes.Item(7)
$s15OrigHeight = $s15Shape.Height
$s15OrigWidth = $s15Shape.Width
$s15Text = $s15Shape.TextFrame.TextRange
$idx2 = $s15Text.Text.IndexOf(" Other?")
$s15Sub = $s15Text.Characters($idx2 + 1, 7)
$s15Sub.Text = " (Topic 8) Scheduling the Trading Bot"
# editing the run can nudge the auto-fit shape size; put it back exactly
$s15Shape.Height = $s15OrigHeight
$s15Shape.Width = $s15OrigWidth
